# "changes for second time"
#
# The login-test sheet stored a real username/password pair in column F
# ("Data") next to the locator rows for the User Name and Password fields.
# This edit wipes those two values (and the mailto: hyperlink that had been
# attached to the password cell), and makes "test_login" the active sheet/
# selection again (it had drifted to "Vessel Master").

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("test_login")

# Row 2 = User Name step -> clear the stored value "raypk".
$wsLogin.Range("F2").ClearContents() | Out-Null

# Row 3 = Password step -> clear the stored value "Homeserver@123".
$wsLogin.Range("F3").ClearContents() | Out-Null

# That password cell carried a mailto: hyperlink (and its Hyperlink cell
# style) - drop the hyperlink itself; the (now blank) cell keeps its style.
$wsLogin.Hyperlinks.Delete() | Out-Null

# Bring "test_login" back to the front, with F2 as the selected cell,
# matching the workbook's saved view state after the edit.
$wsLogin.Activate() | Out-Null
$wsLogin.Range("F2").Select() | Out-Null
